$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9371385197543067
$ws.Range("J2").Value = 0.9371385197543067
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.752937333333333
$ws.Range("N2").Value = 11.258812
$ws.Range("O2").Value = 0.6855621274031838
$ws.Range("P2").Value = 0.6855621274031838
$ws.Range("Q2").Value = 51.77776520523377
$ws.Range("R2").Value = 465.9998868471039
$ws.Range("S2").Value = 0.6424666772742331
$ws.Range("T2").Value = 0.6424666772742331

# Row 3
$ws.Range("I3").Value = 0.9371385197543067
$ws.Range("J3").Value = 0.9371385197543067
$ws.Range("O3").Value = 0.2368266084628361
$ws.Range("P3").Value = 0.2368266084628362
$ws.Range("S3").Value = 0.221939337293295
$ws.Range("T3").Value = 0.2219393372932951

# Row 4
$ws.Range("I4").Value = 0.9371385197543067
$ws.Range("J4").Value = 0.9371385197543067
$ws.Range("M4").Value = 0.4248633333333334
$ws.Range("N4").Value = 1.27459
$ws.Range("O4").Value = 0.07761126413398003
$ws.Range("P4").Value = 0.07761126413398005
$ws.Range("Q4").Value = 5.861668331697778
$ws.Range("R4").Value = 52.75501498528001
$ws.Range("S4").Value = 0.07273250518677857
$ws.Range("T4").Value = 0.07273250518677858

# Row 5
$ws.Range("G5").Value = 0.9254496666666667
$ws.Range("H5").Value = 2.776349
$ws.Range("I5").Value = 0.0628614802456932
$ws.Range("J5").Value = 0.06286148024569319
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.752937333333333
$ws.Range("N5").Value = 11.258812
$ws.Range("O5").Value = 0.6855621274031838
$ws.Range("P5").Value = 0.6855621274031838
$ws.Range("Q5").Value = 3.473154604154222
$ws.Range("R5").Value = 31.258391437388
$ws.Range("S5").Value = 0.04309545012895064
$ws.Range("T5").Value = 0.04309545012895064

# Row 6
$ws.Range("G6").Value = 0.9254496666666667
$ws.Range("H6").Value = 2.776349
$ws.Range("I6").Value = 0.0628614802456932
$ws.Range("J6").Value = 0.06286148024569319
$ws.Range("O6").Value = 0.2368266084628361
$ws.Range("P6").Value = 0.2368266084628362
$ws.Range("Q6").Value = 1.199797060967444
$ws.Range("R6").Value = 10.798173548707
$ws.Range("S6").Value = 0.01488727116954109
$ws.Range("T6").Value = 0.01488727116954109

# Row 7
$ws.Range("G7").Value = 0.9254496666666667
$ws.Range("H7").Value = 2.776349
$ws.Range("I7").Value = 0.0628614802456932
$ws.Range("J7").Value = 0.06286148024569319
$ws.Range("M7").Value = 0.4248633333333334
$ws.Range("N7").Value = 1.27459
$ws.Range("O7").Value = 0.07761126413398003
$ws.Range("P7").Value = 0.07761126413398005
$ws.Range("Q7").Value = 0.3931896302122223
$ws.Range("R7").Value = 3.53870667191
$ws.Range("S7").Value = 0.004878758947201463
$ws.Range("T7").Value = 0.004878758947201463
